$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 38) below the existing table: a text label "7038"
# in column A (stored as a shared string, same pattern as the other labels
# like "7002", "7033", etc.) and a numeric 0 in column B.
#
# The leading apostrophe forces Excel to store the numeric-looking label as
# text rather than silently coercing it to a number; resetting the style
# back to "Normal" afterwards keeps the cell on the default style (matching
# every other label cell in the sheet) instead of leaving it tagged with the
# quote-prefix formatting that the apostrophe entry would otherwise apply.
$ws.Range("A38").Value = "'7038"
$ws.Range("A38").Style = "Normal"
$ws.Range("B38").Value = 0
